# Atualizacao de bases das ligas, do dia: 03-04-2024 as 22:09
# Re-syncs several "Czech Republic First League" match rows whose fixture
# id/odds data had been assigned to the wrong row (a same-date pair got
# swapped), plus a handful of closing-odds corrections for later matches.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 123/124 (match ids 6803057 / 6803053) - swap back to correct rows
$ws.Range("B123").Value = 6803053
$ws.Range("F123").Value = "Hradec Kralove"
$ws.Range("G123").Value = "Viktoria Plzen"
$ws.Range("H123").Value = 1
$ws.Range("I123").Value = 1
$ws.Range("J123").Value = "D"
$ws.Range("K123").Value = 4
$ws.Range("L123").Value = 3.6
$ws.Range("M123").Value = 1.8
$ws.Range("N123").Value = 4.333
$ws.Range("O123").Value = 3.75
$ws.Range("P123").Value = 1.7
$ws.Range("Q123").Value = 0.75
$ws.Range("R123").Value = 1.9
$ws.Range("S123").Value = 1.95
$ws.Range("T123").Value = 2.75
$ws.Range("U123").Value = 2.025
$ws.Range("V123").Value = 1.825
$ws.Range("X123").Value = 2.75
$ws.Range("Y123").Value = -1
$ws.Range("Z123").Value = 0.8999999999999999
$ws.Range("AA123").Value = -1
$ws.Range("AB123").Value = -1
$ws.Range("AC123").Value = 0.825
$ws.Range("B124").Value = 6803057
$ws.Range("F124").Value = "MFK Karvina"
$ws.Range("G124").Value = "Sigma Olomouc"
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 2
$ws.Range("J124").Value = "A"
$ws.Range("K124").Value = 3.1
$ws.Range("L124").Value = 3.2
$ws.Range("M124").Value = 2.25
$ws.Range("N124").Value = 3.4
$ws.Range("O124").Value = 3.1
$ws.Range("P124").Value = 2.15
$ws.Range("Q124").Value = 0.25
$ws.Range("R124").Value = 1.925
$ws.Range("S124").Value = 1.925
$ws.Range("T124").Value = 2.25
$ws.Range("U124").Value = 1.8
$ws.Range("V124").Value = 2.05
$ws.Range("X124").Value = -1
$ws.Range("Y124").Value = 1.15
$ws.Range("Z124").Value = -1
$ws.Range("AA124").Value = 0.925
$ws.Range("AB124").Value = -0.5
$ws.Range("AC124").Value = 0.5249999999999999

# Rows 131/132 (match ids 6803065 / 7521554)
$ws.Range("B131").Value = 6803065
$ws.Range("F131").Value = "FK Jablonec"
$ws.Range("G131").Value = "FK Teplice"
$ws.Range("H131").Value = 3
$ws.Range("I131").Value = 2
$ws.Range("J131").Value = "H"
$ws.Range("K131").Value = 1.75
$ws.Range("L131").Value = 3.4
$ws.Range("M131").Value = 4.75
$ws.Range("N131").Value = 1.6
$ws.Range("O131").Value = 3.75
$ws.Range("P131").Value = 5.5
$ws.Range("Q131").Value = -0.75
$ws.Range("R131").Value = 1.8
$ws.Range("S131").Value = 2.05
$ws.Range("U131").Value = 2.025
$ws.Range("V131").Value = 1.825
$ws.Range("W131").Value = 0.6000000000000001
$ws.Range("Y131").Value = -1
$ws.Range("Z131").Value = 0.4
$ws.Range("AA131").Value = -0.5
$ws.Range("AB131").Value = 1.025
$ws.Range("B132").Value = 7521554
$ws.Range("F132").Value = "Mlada Boleslav"
$ws.Range("G132").Value = "Viktoria Plzen"
$ws.Range("H132").Value = 1
$ws.Range("I132").Value = 3
$ws.Range("J132").Value = "A"
$ws.Range("K132").Value = 3.2
$ws.Range("L132").Value = 3.25
$ws.Range("M132").Value = 2.2
$ws.Range("N132").Value = 3.4
$ws.Range("O132").Value = 3.4
$ws.Range("P132").Value = 2.1
$ws.Range("Q132").Value = 0.25
$ws.Range("R132").Value = 2.025
$ws.Range("S132").Value = 1.825
$ws.Range("U132").Value = 1.825
$ws.Range("V132").Value = 2.025
$ws.Range("W132").Value = -1
$ws.Range("Y132").Value = 1.1
$ws.Range("Z132").Value = -1
$ws.Range("AA132").Value = 0.825
$ws.Range("AB132").Value = 0.825

# Rows 144/145 (match ids 6803080 / 6803081)
$ws.Range("B144").Value = 6803080
$ws.Range("F144").Value = "Mlada Boleslav"
$ws.Range("G144").Value = "Slovacko"
$ws.Range("H144").Value = 0
$ws.Range("I144").Value = 1
$ws.Range("J144").Value = "A"
$ws.Range("K144").Value = 2.3
$ws.Range("L144").Value = 3.3
$ws.Range("M144").Value = 3.1
$ws.Range("N144").Value = 2.5
$ws.Range("O144").Value = 3.3
$ws.Range("P144").Value = 2.8
$ws.Range("Q144").Value = 0
$ws.Range("R144").Value = 1.8
$ws.Range("S144").Value = 2.05
$ws.Range("U144").Value = 1.95
$ws.Range("V144").Value = 1.9
$ws.Range("W144").Value = -1
$ws.Range("Y144").Value = 1.8
$ws.Range("Z144").Value = -1
$ws.Range("AA144").Value = 1.05
$ws.Range("AB144").Value = -1
$ws.Range("AC144").Value = 0.8999999999999999
$ws.Range("B145").Value = 6803081
$ws.Range("F145").Value = "FC Trinity Zlin"
$ws.Range("G145").Value = "Hradec Kralove"
$ws.Range("H145").Value = 4
$ws.Range("I145").Value = 0
$ws.Range("J145").Value = "H"
$ws.Range("K145").Value = 2.5
$ws.Range("L145").Value = 3.4
$ws.Range("M145").Value = 2.7
$ws.Range("N145").Value = 3
$ws.Range("O145").Value = 3.4
$ws.Range("P145").Value = 2.3
$ws.Range("Q145").Value = 0.25
$ws.Range("R145").Value = 1.825
$ws.Range("S145").Value = 2.025
$ws.Range("U145").Value = 2.05
$ws.Range("V145").Value = 1.8
$ws.Range("W145").Value = 2
$ws.Range("Y145").Value = -1
$ws.Range("Z145").Value = 0.825
$ws.Range("AA145").Value = -1
$ws.Range("AB145").Value = 1.05
$ws.Range("AC145").Value = -1

# Rows 191/192 (match ids 6803128 / 6803123)
$ws.Range("B191").Value = 6803128
$ws.Range("F191").Value = "Slavia Prague"
$ws.Range("G191").Value = "FK Teplice"
$ws.Range("H191").Value = 4
$ws.Range("I191").Value = 0
$ws.Range("J191").Value = "H"
$ws.Range("K191").Value = 1.2
$ws.Range("L191").Value = 6.5
$ws.Range("M191").Value = 10
$ws.Range("N191").Value = 1.142
$ws.Range("O191").Value = 6.5
$ws.Range("P191").Value = 17
$ws.Range("Q191").Value = -2
$ws.Range("R191").Value = 1.875
$ws.Range("S191").Value = 1.975
$ws.Range("T191").Value = 3
$ws.Range("U191").Value = 1.975
$ws.Range("V191").Value = 1.875
$ws.Range("W191").Value = 0.1419999999999999
$ws.Range("X191").Value = -1
$ws.Range("Z191").Value = 0.875
$ws.Range("AA191").Value = -1
$ws.Range("AB191").Value = 0.9750000000000001
$ws.Range("B192").Value = 6803123
$ws.Range("F192").Value = "Ceske Budejovice"
$ws.Range("G192").Value = "FC Trinity Zlin"
$ws.Range("H192").Value = 2
$ws.Range("I192").Value = 2
$ws.Range("J192").Value = "D"
$ws.Range("K192").Value = 1.833
$ws.Range("L192").Value = 3.6
$ws.Range("M192").Value = 3.8
$ws.Range("N192").Value = 1.909
$ws.Range("O192").Value = 3.5
$ws.Range("P192").Value = 3.6
$ws.Range("Q192").Value = -0.5
$ws.Range("R192").Value = 2
$ws.Range("S192").Value = 1.85
$ws.Range("T192").Value = 2.5
$ws.Range("U192").Value = 1.85
$ws.Range("V192").Value = 2
$ws.Range("W192").Value = -1
$ws.Range("X192").Value = 2.5
$ws.Range("Z192").Value = -1
$ws.Range("AA192").Value = 0.8500000000000001
$ws.Range("AB192").Value = 0.8500000000000001

# Rows 202/203 (match ids 6803139 / 6803142)
$ws.Range("B202").Value = 6803139
$ws.Range("F202").Value = "Bohemians 1905"
$ws.Range("G202").Value = "FK Jablonec"
$ws.Range("H202").Value = 2
$ws.Range("I202").Value = 0
$ws.Range("J202").Value = "H"
$ws.Range("K202").Value = 2.1
$ws.Range("L202").Value = 3.4
$ws.Range("M202").Value = 3.5
$ws.Range("N202").Value = 2.25
$ws.Range("O202").Value = 3.25
$ws.Range("P202").Value = 3.3
$ws.Range("Q202").Value = -0.25
$ws.Range("R202").Value = 1.85
$ws.Range("S202").Value = 2
$ws.Range("T202").Value = 2.5
$ws.Range("U202").Value = 1.975
$ws.Range("V202").Value = 1.875
$ws.Range("W202").Value = 1.25
$ws.Range("Y202").Value = -1
$ws.Range("Z202").Value = 0.8500000000000001
$ws.Range("AA202").Value = -1
$ws.Range("AB202").Value = -1
$ws.Range("AC202").Value = 0.875
$ws.Range("B203").Value = 6803142
$ws.Range("F203").Value = "Viktoria Plzen"
$ws.Range("G203").Value = "Slovan Liberec"
$ws.Range("H203").Value = 1
$ws.Range("I203").Value = 3
$ws.Range("J203").Value = "A"
$ws.Range("K203").Value = 1.4
$ws.Range("L203").Value = 5
$ws.Range("M203").Value = 7
$ws.Range("N203").Value = 1.4
$ws.Range("O203").Value = 4.75
$ws.Range("P203").Value = 6.5
$ws.Range("Q203").Value = -1.25
$ws.Range("R203").Value = 1.975
$ws.Range("S203").Value = 1.875
$ws.Range("T203").Value = 2.75
$ws.Range("U203").Value = 1.9
$ws.Range("V203").Value = 1.95
$ws.Range("W203").Value = -1
$ws.Range("Y203").Value = 5.5
$ws.Range("Z203").Value = -1
$ws.Range("AA203").Value = 0.875
$ws.Range("AB203").Value = 0.8999999999999999
$ws.Range("AC203").Value = -1

# Rows 210-217 - closing odds corrections (no id/team swap)
$ws.Range("R210").Value = 2.05
$ws.Range("S210").Value = 1.8
$ws.Range("N211").Value = 2.45
$ws.Range("P211").Value = 2.875
$ws.Range("R211").Value = 1.8
$ws.Range("S211").Value = 2.05
$ws.Range("N212").Value = 2.45
$ws.Range("O212").Value = 2.9
$ws.Range("P212").Value = 3.2
$ws.Range("R212").Value = 2.1
$ws.Range("S212").Value = 1.775
$ws.Range("U212").Value = 1.875
$ws.Range("V212").Value = 1.975
$ws.Range("R213").Value = 1.9
$ws.Range("S213").Value = 1.95
$ws.Range("U214").Value = 1.875
$ws.Range("V214").Value = 1.975
$ws.Range("U215").Value = 1.975
$ws.Range("V215").Value = 1.875
$ws.Range("R216").Value = 1.975
$ws.Range("S216").Value = 1.875
$ws.Range("U216").Value = 2.025
$ws.Range("V216").Value = 2
$ws.Range("T217").Value = 3
$ws.Range("U217").Value = 1.825
$ws.Range("V217").Value = 2.025
